$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "SA"
$ws.Range("L13").Value = 0.9938510866456975

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.007981585008737
$ws.Range("D16").Value = 0.9702898268664496
$ws.Range("E16").Value = 0.998104562345598
$ws.Range("F16").Value = 0.9911618465239928
$ws.Range("G16").Value = 1.007981585008737
$ws.Range("H16").Value = 0.9702898268664496
$ws.Range("I16").Value = 1.004081529301284
$ws.Range("J16").Value = 0.9884821335975028
$ws.Range("K16").Value = 0.9999638822424606
$ws.Range("L16").Value = 0.9791718525351845
$ws.Range("M16").Value = 1.007981585008737
$ws.Range("N16").Value = 0.9841971946060237
$ws.Range("O16").Value = 0.9918844551861942
$ws.Range("P16").Value = 0.9924046523026512

$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14
$excel.CutCopyMode = $false
